$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 87.442999999999998
$ws.Range("D2").Value = 0.055899999999999998

$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 0.60199999999999998
$ws.Range("D5").Value = 0.055899999999999998

$ws.Range("K5").Select()
